$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 4075
$ws.Range("L3").Value = 4305
$ws.Range("L4").Value = 1056
$ws.Range("L5").Value = 247
$ws.Range("L6").Value = 3702
$ws.Range("L7").Value = 13385

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L3").Value = 41
$ws.Range("L7").Value = 154

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 256
$ws.Range("L3").Value = 296
$ws.Range("L6").Value = 240
$ws.Range("L7").Value = 882

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L3").Value = 209
$ws.Range("L6").Value = 197
$ws.Range("L7").Value = 628

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L3").Value = 58
$ws.Range("L7").Value = 183

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L5").Value = 17
$ws.Range("L7").Value = 488

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L3").Value = 77
$ws.Range("L7").Value = 254

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L3").Value = 97
$ws.Range("L7").Value = 222

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L7").Value = 442
$ws.Range("L8").Value = 882
$ws.Range("L9").Value = 80
$ws.Range("L11").Value = 216
$ws.Range("L15").Value = 98
$ws.Range("L19").Value = 381
$ws.Range("L20").Value = 335
$ws.Range("L24").Value = 33
$ws.Range("L29").Value = 752
$ws.Range("L31").Value = 129
$ws.Range("L33").Value = 628
$ws.Range("L34").Value = 81
$ws.Range("L37").Value = 488
$ws.Range("L41").Value = 58
$ws.Range("L42").Value = 426
$ws.Range("L43").Value = 103
$ws.Range("L48").Value = 177
$ws.Range("L51").Value = 164
$ws.Range("L53").Value = 154
$ws.Range("L54").Value = 278
$ws.Range("L55").Value = 127
$ws.Range("L56").Value = 11
$ws.Range("L63").Value = 41
$ws.Range("L65").Value = 254
$ws.Range("L68").Value = 41
$ws.Range("L72").Value = 57
$ws.Range("L75").Value = 48
$ws.Range("L78").Value = 168
$ws.Range("L79").Value = 351
$ws.Range("L85").Value = 704
$ws.Range("L86").Value = 99
$ws.Range("L90").Value = 133
$ws.Range("L91").Value = 188
$ws.Range("L95").Value = 183
$ws.Range("L96").Value = 142
$ws.Range("L99").Value = 222
$ws.Range("L101").Value = 13385

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L2").Value = 46
$ws.Range("L7").Value = 129

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L6").Value = 134
$ws.Range("L7").Value = 278

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L3").Value = 283
$ws.Range("L6").Value = 192
$ws.Range("L7").Value = 752

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L6").Value = 77
$ws.Range("L7").Value = 177

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L3").Value = 117
$ws.Range("L7").Value = 381

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L3").Value = 138
$ws.Range("L6").Value = 121
$ws.Range("L7").Value = 426

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L6").Value = 51
$ws.Range("L7").Value = 168

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L3").Value = 44
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 33

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L2").Value = 47
$ws.Range("L7").Value = 142

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L2").Value = 67
$ws.Range("L3").Value = 78
$ws.Range("L7").Value = 188

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 117
$ws.Range("L7").Value = 351

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 106
$ws.Range("L7").Value = 335

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L3").Value = 138
$ws.Range("L7").Value = 442

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L3").Value = 31
$ws.Range("L7").Value = 98

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L3").Value = 65
$ws.Range("L7").Value = 216

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 80

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 99

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L3").Value = 43
$ws.Range("L7").Value = 133

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("L2").Value = 13
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L3").Value = 29
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 215
$ws.Range("L3").Value = 283
$ws.Range("L5").Value = 14
$ws.Range("L7").Value = 704

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("L2").Value = 3
$ws.Range("L7").Value = 11
